# edit.ps1 - applies the "Last changes before presentation" revision
#
# Four edits to the document body text:
#  1. "fondamental du barometre"   -> "fondamental d'un outil bien connu, le barometre"
#     (splits the run into 3 runs at the inserted-text boundaries)
#  2. "dans le cours a la piezo"   -> "dans le cours sur les capteurs a la piezo"
#     (splits the run into 3 runs, and the "_GoBack" bookmark is relocated into
#      the new split point, between " sur les capteurs" and " a la piezo...")
#  3. "mecanique. Donc on utilise" -> "mecanique. On utilise"
#     (splits the run into 3 runs: "...mecanique. " | "O" | "n utilise...", while
#      keeping the other, untouched runs later in that same paragraph intact)
#  4. The "_GoBack" bookmark that used to sit right after "DEMO :" is removed
#     (it is the same bookmark that gets moved in edit #2 - Word bookmark names
#      are unique, so moving it removes it from its old spot automatically)

$d = $word.ActiveDocument

function Split-RangeFormatting($rng) {
    # Forces Word to break the run(s) covering $rng away from their neighbours
    # by nudging a formatting property and setting it straight back - this
    # mirrors what Word itself does internally when a sub-range ends up with
    # distinct run properties, without altering the rendered appearance.
    $sz = $rng.Font.Size
    $rng.Font.Size = $sz + 1
    $rng.Font.Size = $sz
}

function Find-Range($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, `
        $true, 1, $false, "", 0) | Out-Null
    return $r
}

# ---------------------------------------------------------------------------
# Change 1: "... fondamental du barometre ..." ->
#           "... fondamental d'un outil bien connu, le barometre ..."
# (this paragraph only ever contains a single run, so no neighbouring runs
# are at risk of being merged together by the text edit)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("fondamental du baromètre", $true, $false, $false, $false, $false, `
    $true, 1, $false, "fondamental d'un outil bien connu, le baromètre", 2) | Out-Null

$newPhrase = Find-Range("d'un outil bien connu, le")
Split-RangeFormatting $newPhrase

# ---------------------------------------------------------------------------
# Change 2: "... dans le cours a la piezo-electricite ..." ->
#           "... dans le cours sur les capteurs [_GoBack] a la piezo-electricite ..."
# (this paragraph, too, only ever contains a single run)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("dans le cours à la piézo", $true, $false, $false, $false, $false, `
    $true, 1, $false, "dans le cours sur les capteurs à la piézo", 2) | Out-Null

$newPhrase2 = Find-Range(" sur les capteurs")
Split-RangeFormatting $newPhrase2

# Relocate the "_GoBack" bookmark to sit right after " sur les capteurs"
# (it currently exists right after "DEMO :" further down the document).
$goBackTarget = $d.Range($newPhrase2.End, $newPhrase2.End)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $goBackTarget) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "... mecanique. Donc on utilise ..." -> "... mecanique. On utilise ..."
# This paragraph has several sibling runs after the edited text ("des cristaux"
# + " " + "dielectriques..." + " " + "electriques...mecaniques." + " "). A
# text-length-changing edit in this runtime re-merges every identically
# formatted run in the paragraph, so after the replace we explicitly restore
# those original run boundaries (in addition to the two brand-new ones).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("mécanique. Donc on utilise", $true, $false, $false, $false, $false, `
    $true, 1, $false, "mécanique. On utilise", 2) | Out-Null

# New split: "...mecanique. " | "O" | "n utilise..."
$oFull = Find-Range("mécanique. On utilise")
$oStart = $oFull.Start + 11   # length of "mécanique. " (incl. trailing space)
$oOnly = $d.Range($oStart, $oStart + 1)
Split-RangeFormatting $oOnly

# Restore the original (untouched) run boundaries further along the paragraph
$space1 = Find-Range(" diélectriques")
Split-RangeFormatting ($d.Range($space1.Start, $space1.Start + 1))

Split-RangeFormatting (Find-Range("diélectriques comme le quartz cristallin qui font apparaître des charges "))

Split-RangeFormatting (Find-Range("électriques lorsqu’ils sont soumis à des contraintes mécaniques."))

# ---------------------------------------------------------------------------
# Change 4: the old "_GoBack" bookmark after "DEMO :" no longer exists - this
# was already achieved above via Delete()+Add() in Change 2.
# ---------------------------------------------------------------------------

Write-Output "Done."
